$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; this shifts existing rows 40..151 down to 41..152
# (Excel automatically copies formatting, including the date-style on column D, from
# the row being pushed down, so the new row 40 already carries the right styles.)
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new weekly price observation.
$ws.Cells.Item(40, 1).Value = 8
$ws.Cells.Item(40, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(40, 3).Value = "Coquimbo"
$ws.Cells.Item(40, 4).Value = 44414
$ws.Cells.Item(40, 5).Value = 4
$ws.Cells.Item(40, 6).Value = 100112032
$ws.Cells.Item(40, 7).Value = "Zapallo italiano"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 600
$ws.Cells.Item(40, 11).Value = 7500
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = 7750
$ws.Cells.Item(40, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(40, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value = 155
$ws.Cells.Item(40, 17).Value = 50
$ws.Cells.Item(40, 18).Value = "Hortaliza"
